$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Stimulus"

# Rename the stimulus test case text. The old string
# "Multiple Write (Any) (3 - 10 txns) + ..." is no longer referenced by any
# cell once this is applied, so it drops out of the shared-string table and
# everything after it reindexes down by one; the new (longer) string is
# appended at the end of the table - exactly the shape of the target diff.
$ws1.Range("C6").Value = "Non-consecutive Multiple Write (Any) (3 - 10 txns) + Multiple Read (Any) (3 - 10 txns) (Only trigger read after write)"

# Widen column C on the Stimulus sheet so the longer text keeps fitting
# (closest reachable width to the authored 104.54296875).
$ws1.Columns.Item(3).ColumnWidth = 103.73

# Make "Stimulus" the active/selected sheet again (it was "Coverage"
# before), and leave the last-used selection on each sheet matching the
# edit: E20 on Stimulus, F22 (unchanged) on Coverage.
$ws1.Activate()
$ws1.Range("E20").Select() | Out-Null
